# Fix a couple of typos in the ENERGY_CARRIERS database and update the
# active cell selection, as described in the commit:
# "Adding 'ENERGY_CARRIERS' to reference-case - Fixing a few typos in the
#  database and adding it to the reference-case to fix workflows."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ENERGY_CARRIERS")

# Fix typos in the "description" column (column A)
# Row 17: "Sulight" -> "Sunlight"
$ws.Range("A17").Value = "Sunlight"
# Row 18: "ultraviolet" -> "Ultraviolet"
$ws.Range("A18").Value = "Ultraviolet"

# Update the selected cell in the sheet view to A27
$ws.Range("A27").Select()
